$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio3")

# ---------------------------------------------------------------------------
# Row 45 - header row: "test" -> "test 1", add new mirrored table at F:I
# ---------------------------------------------------------------------------
$ws.Range("A45").Value = "test 1"

$ws.Range("A45:D45").Copy()
$ws.Range("F45").PasteSpecial(-4122)
$ws.Range("F45").Value = "test 2"
$ws.Range("G45").Value = "cal"
$ws.Range("H45").Value = "pro"
$ws.Range("I45").Value = "eur"

# ---------------------------------------------------------------------------
# Row 46 - 70g prot whey no g (mirror into F:I)
# ---------------------------------------------------------------------------
$ws.Range("A46:D46").Copy()
$ws.Range("F46").PasteSpecial(-4122)
$ws.Range("F46").Value = "70g prot whey no g"
$ws.Range("G46").Formula = "=`$O`$12*0.7"
$ws.Range("H46").Formula = "=`$P`$12*0.7"
$ws.Range("I46").Formula = "=`$Q`$12*0.7"

# ---------------------------------------------------------------------------
# Row 47 - 500ml latte esse: fix F47 (was "sistemare", unstyled) + mirror G:I
# ---------------------------------------------------------------------------
$ws.Range("A47:D47").Copy()
$ws.Range("F47").PasteSpecial(-4122)
$ws.Range("F47").Value = "500ml latte esse"
$ws.Range("G47").Formula = "=47*5"
$ws.Range("H47").Formula = "=3.3*5"
$ws.Range("I47").Formula = "=0.1*5"

# ---------------------------------------------------------------------------
# Row 48 - 2 wurstel smart maxi (mirror into F:I)
# ---------------------------------------------------------------------------
$ws.Range("A48:D48").Copy()
$ws.Range("F48").PasteSpecial(-4122)
$ws.Range("F48").Value = "2 wurstel smart maxi"
$ws.Range("G48").Formula = "=240*2"
$ws.Range("H48").Formula = "=15*2"
$ws.Range("I48").Formula = "=0.26*2"

# ---------------------------------------------------------------------------
# Row 49 - 100g pane (mirror into F:I)
# ---------------------------------------------------------------------------
$ws.Range("A49:D49").Copy()
$ws.Range("F49").PasteSpecial(-4122)
$ws.Range("F49").Value = "100g pane"
$ws.Range("G49").Formula = "=`$G`$6*1"
$ws.Range("H49").Formula = "=`$H`$6*1"
$ws.Range("I49").Formula = "=`$I`$6*1"

# ---------------------------------------------------------------------------
# Row 50 - 500g yogurt bianco s: quantity 5 -> 2.5, D50 style 5 -> 19, mirror F:I
# ---------------------------------------------------------------------------
$ws.Range("B50").Formula = "=51*2.5"
$ws.Range("C50").Formula = "=5.3*2.5"
$ws.Range("D49").Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("D50").Formula = "=0.15*2.5"

$ws.Range("A50:D50").Copy()
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("F50").ClearContents()
$ws.Range("G50").ClearContents()
$ws.Range("H50").ClearContents()
$ws.Range("I50").ClearContents()

# ---------------------------------------------------------------------------
# Row 51 - empty spacer row, mirror formats into F:I
# ---------------------------------------------------------------------------
$ws.Range("A51:D51").Copy()
$ws.Range("F51").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 52 - empty spacer row, mirror formats into F:I
# ---------------------------------------------------------------------------
$ws.Range("A52:D52").Copy()
$ws.Range("F52").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 53 - totale row, mirror sums into F:I
# ---------------------------------------------------------------------------
$ws.Range("A53:D53").Copy()
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("F53").Value = "totale"
$ws.Range("G53").Formula = "=SUM(G46:G52)"
$ws.Range("H53").Formula = "=SUM(H46:H52)"
$ws.Range("I53").Formula = "=SUM(I46:I52)"

# ---------------------------------------------------------------------------
# Row 55 - header row: "cena" -> "cena 1", mirror into F:I
# ---------------------------------------------------------------------------
$ws.Range("A55").Value = "cena 1"

$ws.Range("A55:D55").Copy()
$ws.Range("F55").PasteSpecial(-4122)
$ws.Range("F55").Value = "cena 1"
$ws.Range("G55").Value = "cal"
$ws.Range("H55").Value = "pro"
$ws.Range("I55").Value = "eur"

# ---------------------------------------------------------------------------
# Row 56 - 1 wurstel smart maxi, mirror new "35g prot whey no g" into F:I
# ---------------------------------------------------------------------------
$ws.Range("A56:D56").Copy()
$ws.Range("F56").PasteSpecial(-4122)
$ws.Range("F56").Value = "35g prot whey no g"
$ws.Range("G56").Formula = "=`$O`$12*0.35"
$ws.Range("H56").Formula = "=`$P`$12*0.35"
$ws.Range("I56").Formula = "=`$Q`$12*0.35"

# ---------------------------------------------------------------------------
# Row 57 - "50g pane" -> "100g pasta smart", formulas -> constants; mirror F:I
# ---------------------------------------------------------------------------
$ws.Range("A57").Value = "100g pasta smart"
$ws.Range("B57").Value = 350
$ws.Range("C57").Value = 11.5
$ws.Range("D49").Copy()
$ws.Range("D57").PasteSpecial(-4122)
$ws.Range("D57").Value = 0.13

$ws.Range("A57:D57").Copy()
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("F57").Value = "250ml latte esse"
$ws.Range("G57").Formula = "=47*2.5"
$ws.Range("H57").Formula = "=3.3*2.5"
$ws.Range("I57").Formula = "=0.1*2.5"

# ---------------------------------------------------------------------------
# Row 58 - old "tot" row content cleared (kept as blank styled cells);
#          mirror "100g pasta smart" (constants) into F:I
# ---------------------------------------------------------------------------
$ws.Range("A58:D58").Copy()
$ws.Range("F58").PasteSpecial(-4122)

$ws.Range("A58").ClearContents()
$ws.Range("B58").ClearContents()
$ws.Range("C58").ClearContents()
$ws.Range("D58").ClearContents()

$ws.Range("F58").Value = "100g pasta smart"
$ws.Range("G58").Value = 350
$ws.Range("H58").Value = 11.5
$ws.Range("I58").Value = 0.13

# ---------------------------------------------------------------------------
# Row 59 (NEW) - "tot" row summing rows 56:58, mirror into F:I
# ---------------------------------------------------------------------------
$ws.Range("A53").Copy()
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("A59").Value = "tot"

$ws.Range("B53").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$ws.Range("B59").Formula = "=B53+SUM(B56:B58)"

$ws.Range("C53").Copy()
$ws.Range("C59").PasteSpecial(-4122)
$ws.Range("C59").Formula = "=C53+SUM(C56:C58)"

$ws.Range("D46").Copy()
$ws.Range("D59").PasteSpecial(-4122)
$ws.Range("D59").Formula = "=D53+SUM(D56:D58)"

$ws.Range("A59:D59").Copy()
$ws.Range("F59").PasteSpecial(-4122)
$ws.Range("F59").Value = "tot"
$ws.Range("G59").Formula = "=G53+SUM(G56:G58)"
$ws.Range("H59").Formula = "=H53+SUM(H56:H58)"
$ws.Range("I59").Formula = "=I53+SUM(I56:I58)"

# ---------------------------------------------------------------------------
# Row 61 (NEW) - header row "temporaneo", mirror into F:I
# ---------------------------------------------------------------------------
$ws.Range("A55:D55").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("A61").Value = "temporaneo"
$ws.Range("B61").Value = "cal"
$ws.Range("C61").Value = "pro"
$ws.Range("D61").Value = "eur"

$ws.Range("A61:D61").Copy()
$ws.Range("F61").PasteSpecial(-4122)
$ws.Range("F61").Value = "temporaneo"
$ws.Range("G61").Value = "cal"
$ws.Range("H61").Value = "pro"
$ws.Range("I61").Value = "eur"

# ---------------------------------------------------------------------------
# Row 62 (NEW) - 2 wurstel smart maxi, mirror new "35g prot whey no g" F:I
# ---------------------------------------------------------------------------
$ws.Range("A48:D48").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A62").Value = "2 wurstel smart maxi"
$ws.Range("B62").Formula = "=240*2"
$ws.Range("C62").Formula = "=15*2"
$ws.Range("D62").Formula = "=0.26*2"

$ws.Range("A62:D62").Copy()
$ws.Range("F62").PasteSpecial(-4122)
$ws.Range("F62").Value = "35g prot whey no g"
$ws.Range("G62").Formula = "=`$O`$12*0.35"
$ws.Range("H62").Formula = "=`$P`$12*0.35"
$ws.Range("I62").Formula = "=`$Q`$12*0.35"

# ---------------------------------------------------------------------------
# Row 63 (NEW) - empty spacer row, mirror "250ml latte esse" into F:I
# ---------------------------------------------------------------------------
$ws.Range("A51:D51").Copy()
$ws.Range("A63").PasteSpecial(-4122)

$ws.Range("A63:D63").Copy()
$ws.Range("F63").PasteSpecial(-4122)
$ws.Range("F63").Value = "250ml latte esse"
$ws.Range("G63").Formula = "=47*2.5"
$ws.Range("H63").Formula = "=3.3*2.5"
$ws.Range("I63").Formula = "=0.1*2.5"

# ---------------------------------------------------------------------------
# Row 64 (NEW) - empty spacer row (A:D), "200g gnocchi smart" constants F:I
# ---------------------------------------------------------------------------
$ws.Range("A58:D58").Copy()
$ws.Range("A64").PasteSpecial(-4122)

$ws.Range("A64:D64").Copy()
$ws.Range("F64").PasteSpecial(-4122)
$ws.Range("F64").Value = "200g gnocchi smart"
$ws.Range("G64").Value = 348
$ws.Range("H64").Value = 9
$ws.Range("I64").Value = 0.2

# ---------------------------------------------------------------------------
# Row 65 (NEW) - "tot" row summing rows 62:64, mirror into F:I
# ---------------------------------------------------------------------------
$ws.Range("A59").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A65").Value = "tot"

$ws.Range("B59").Copy()
$ws.Range("B65").PasteSpecial(-4122)
$ws.Range("B65").Formula = "=B53+SUM(B62:B64)"

$ws.Range("C59").Copy()
$ws.Range("C65").PasteSpecial(-4122)
$ws.Range("C65").Formula = "=C53+SUM(C62:C64)"

$ws.Range("D59").Copy()
$ws.Range("D65").PasteSpecial(-4122)
$ws.Range("D65").Formula = "=D53+SUM(D62:D64)"
$ws.Range("D65").NumberFormat = $ws.Range("D46").NumberFormat

$ws.Range("A65:D65").Copy()
$ws.Range("F65").PasteSpecial(-4122)
$ws.Range("F65").Value = "tot"
$ws.Range("G65").Formula = "=G53+SUM(G62:G64)"
$ws.Range("H65").Formula = "=H53+SUM(H62:H64)"
$ws.Range("I65").Formula = "=I53+SUM(I62:I64)"
$ws.Range("I65").NumberFormat = $ws.Range("D46").NumberFormat

# ---------------------------------------------------------------------------
# Sheet view bookkeeping: scroll position + active selection
# ---------------------------------------------------------------------------
$ws.Range("F46").Select()
$excel.ActiveWindow.ScrollRow = 42

Write-Host "edit complete"
